# Auto-generated Excel COM-interop script
# Implements: insert "2022-Q1" detail sheet before "总计"; prepend 2022-Q1 summary row to "总计"

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# ---- 1. Update 总计 sheet: prepend a 2022-Q1 summary row, shifting existing rows down ----
$totalSheet.Range("A6").Copy($totalSheet.Range("A7"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = '2022-Q1'
$totalSheet.Range("C2").Value = 40
$totalSheet.Range("D2").Value = 20.4

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = '2021-Q4'
$totalSheet.Range("C3").Value = 35
$totalSheet.Range("D3").Value = 20.13

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = '2021-Q3'
$totalSheet.Range("C4").Value = 18
$totalSheet.Range("D4").Value = 12.55

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = '2021-Q2'
$totalSheet.Range("C5").Value = 7
$totalSheet.Range("D5").Value = 5.62

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = '2021-Q1'
$totalSheet.Range("C6").Value = 5
$totalSheet.Range("D6").Value = 5.58

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = '2020-Q4'
$totalSheet.Range("C7").Value = 6
$totalSheet.Range("D7").Value = 5.33

# ---- 2. Insert new "2022-Q1" sheet before "总计" with fund detail data ----
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

$src.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$newSheet.Range("B1").Value = '基金代码'
$newSheet.Range("C1").Value = '基金名称'
$newSheet.Range("D1").Value = '基金规模'
$newSheet.Range("E1").Value = '股票总仓位'
$newSheet.Range("F1").Value = '仓位占比'
$newSheet.Range("G1").Value = '持有市值(亿元)'
$newSheet.Range("H1").Value = '仓位排名'

$src.Range("A2:G2").Copy($newSheet.Range("A2:G2"))
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet "B2" '003751'
$newSheet.Range("C2").Value = '万家瑞隆混合'
Set-TextValue $newSheet "D2" '27.84'
Set-TextValue $newSheet "E2" '86.40'
Set-TextValue $newSheet "F2" '8.38'
Set-TextValue $newSheet "G2" '2.3330'
$newSheet.Range("H2").Value = 1

$src.Range("A2:G2").Copy($newSheet.Range("A3:G3"))
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet "B3" '003634'
$newSheet.Range("C3").Value = '嘉实农业产业股票'
Set-TextValue $newSheet "D3" '28.18'
Set-TextValue $newSheet "E3" '94.29'
Set-TextValue $newSheet "F3" '7.73'
Set-TextValue $newSheet "G3" '2.1783'
$newSheet.Range("H3").Value = 5

$src.Range("A2:G2").Copy($newSheet.Range("A4:G4"))
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet "B4" '519196'
$newSheet.Range("C4").Value = '万家新兴蓝筹灵活配置混合'
Set-TextValue $newSheet "D4" '21.26'
Set-TextValue $newSheet "E4" '80.70'
Set-TextValue $newSheet "F4" '9.44'
Set-TextValue $newSheet "G4" '2.0069'
$newSheet.Range("H4").Value = 1

$src.Range("A2:G2").Copy($newSheet.Range("A5:G5"))
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet "B5" '161810'
$newSheet.Range("C5").Value = '银华内需精选混合(LOF)'
Set-TextValue $newSheet "D5" '25.59'
Set-TextValue $newSheet "E5" '94.71'
Set-TextValue $newSheet "F5" '6.90'
Set-TextValue $newSheet "G5" '1.7657'
$newSheet.Range("H5").Value = 6

$src.Range("A2:G2").Copy($newSheet.Range("A6:G6"))
$newSheet.Range("A6").Value = 4
Set-TextValue $newSheet "B6" '519195'
$newSheet.Range("C6").Value = '万家品质生活灵活配置混合'
Set-TextValue $newSheet "D6" '17.66'
Set-TextValue $newSheet "E6" '79.87'
Set-TextValue $newSheet "F6" '9.51'
Set-TextValue $newSheet "G6" '1.6795'
$newSheet.Range("H6").Value = 2

$src.Range("A2:G2").Copy($newSheet.Range("A7:G7"))
$newSheet.Range("A7").Value = 5
Set-TextValue $newSheet "B7" '009199'
$newSheet.Range("C7").Value = '万家价值优势一年持有期混合'
Set-TextValue $newSheet "D7" '13.70'
Set-TextValue $newSheet "E7" '89.53'
Set-TextValue $newSheet "F7" '9.53'
Set-TextValue $newSheet "G7" '1.3056'
$newSheet.Range("H7").Value = 2

$src.Range("A2:G2").Copy($newSheet.Range("A8:G8"))
$newSheet.Range("A8").Value = 6
Set-TextValue $newSheet "B8" '161912'
$newSheet.Range("C8").Value = '万家社会责任18个月定期开放混合（LOF）A'
Set-TextValue $newSheet "D8" '13.56'
Set-TextValue $newSheet "E8" '88.11'
Set-TextValue $newSheet "F8" '9.57'
Set-TextValue $newSheet "G8" '1.2977'
$newSheet.Range("H8").Value = 1

$src.Range("A2:G2").Copy($newSheet.Range("A9:G9"))
$newSheet.Range("A9").Value = 7
Set-TextValue $newSheet "B9" '005094'
$newSheet.Range("C9").Value = '万家臻选混合'
Set-TextValue $newSheet "D9" '13.43'
Set-TextValue $newSheet "E9" '73.12'
Set-TextValue $newSheet "F9" '8.68'
Set-TextValue $newSheet "G9" '1.1657'
$newSheet.Range("H9").Value = 1

$src.Range("A2:G2").Copy($newSheet.Range("A10:G10"))
$newSheet.Range("A10").Value = 8
Set-TextValue $newSheet "B10" '519181'
$newSheet.Range("C10").Value = '万家和谐增长混合'
Set-TextValue $newSheet "D10" '11.51'
Set-TextValue $newSheet "E10" '88.62'
Set-TextValue $newSheet "F10" '9.43'
Set-TextValue $newSheet "G10" '1.0854'
$newSheet.Range("H10").Value = 2

$src.Range("A2:G2").Copy($newSheet.Range("A11:G11"))
$newSheet.Range("A11").Value = 9
Set-TextValue $newSheet "B11" '005106'
$newSheet.Range("C11").Value = '银华农业产业股票'
Set-TextValue $newSheet "D11" '13.24'
Set-TextValue $newSheet "E11" '93.41'
Set-TextValue $newSheet "F11" '6.28'
Set-TextValue $newSheet "G11" '0.8315'
$newSheet.Range("H11").Value = 5

$src.Range("A2:G2").Copy($newSheet.Range("A12:G12"))
$newSheet.Range("A12").Value = 10
Set-TextValue $newSheet "B12" '161838'
$newSheet.Range("C12").Value = '银华创业板两年定期开放混合'
Set-TextValue $newSheet "D12" '10.44'
Set-TextValue $newSheet "E12" '95.40'
Set-TextValue $newSheet "F12" '5.36'
Set-TextValue $newSheet "G12" '0.5596'
$newSheet.Range("H12").Value = 9

$src.Range("A2:G2").Copy($newSheet.Range("A13:G13"))
$newSheet.Range("A13").Value = 11
Set-TextValue $newSheet "B13" '501083'
$newSheet.Range("C13").Value = '银华科创主题 3 年封闭运作灵活配置混合型证券投资'
Set-TextValue $newSheet "D13" '15.88'
Set-TextValue $newSheet "E13" '59.30'
Set-TextValue $newSheet "F13" '2.78'
Set-TextValue $newSheet "G13" '0.4415'
$newSheet.Range("H13").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A14:G14"))
$newSheet.Range("A14").Value = 12
Set-TextValue $newSheet "B14" '009958'
$newSheet.Range("C14").Value = '长安鑫悦消费驱动混合A'
Set-TextValue $newSheet "D14" '6.98'
Set-TextValue $newSheet "E14" '86.78'
Set-TextValue $newSheet "F14" '6.15'
Set-TextValue $newSheet "G14" '0.4293'
$newSheet.Range("H14").Value = 3

$src.Range("A2:G2").Copy($newSheet.Range("A15:G15"))
$newSheet.Range("A15").Value = 13
Set-TextValue $newSheet "B15" '011817'
$newSheet.Range("C15").Value = '银华阿尔法混合型证券投资基金'
Set-TextValue $newSheet "D15" '10.82'
Set-TextValue $newSheet "E15" '67.94'
Set-TextValue $newSheet "F15" '3.19'
Set-TextValue $newSheet "G15" '0.3452'
$newSheet.Range("H15").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A16:G16"))
$newSheet.Range("A16").Value = 14
Set-TextValue $newSheet "B16" '011733'
$newSheet.Range("C16").Value = '银华瑞祥一年持有期混合型证券投资基金'
Set-TextValue $newSheet "D16" '6.18'
Set-TextValue $newSheet "E16" '72.47'
Set-TextValue $newSheet "F16" '4.69'
Set-TextValue $newSheet "G16" '0.2898'
$newSheet.Range("H16").Value = 5

$src.Range("A2:G2").Copy($newSheet.Range("A17:G17"))
$newSheet.Range("A17").Value = 15
Set-TextValue $newSheet "B17" '161818'
$newSheet.Range("C17").Value = '银华消费主题混合'
Set-TextValue $newSheet "D17" '5.82'
Set-TextValue $newSheet "E17" '90.60'
Set-TextValue $newSheet "F17" '4.61'
Set-TextValue $newSheet "G17" '0.2683'
$newSheet.Range("H17").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A18:G18"))
$newSheet.Range("A18").Value = 16
Set-TextValue $newSheet "B18" '180001'
$newSheet.Range("C18").Value = '银华优势企业混合'
Set-TextValue $newSheet "D18" '6.66'
Set-TextValue $newSheet "E18" '67.80'
Set-TextValue $newSheet "F18" '3.98'
Set-TextValue $newSheet "G18" '0.2651'
$newSheet.Range("H18").Value = 9

$src.Range("A2:G2").Copy($newSheet.Range("A19:G19"))
$newSheet.Range("A19").Value = 17
Set-TextValue $newSheet "B19" '001163'
$newSheet.Range("C19").Value = '银华中国梦30股票'
Set-TextValue $newSheet "D19" '5.00'
Set-TextValue $newSheet "E19" '90.33'
Set-TextValue $newSheet "F19" '5.10'
Set-TextValue $newSheet "G19" '0.2550'
$newSheet.Range("H19").Value = 6

$src.Range("A2:G2").Copy($newSheet.Range("A20:G20"))
$newSheet.Range("A20").Value = 18
Set-TextValue $newSheet "B20" '005481'
$newSheet.Range("C20").Value = '银华瑞泰灵活配置混合'
Set-TextValue $newSheet "D20" '4.80'
Set-TextValue $newSheet "E20" '89.96'
Set-TextValue $newSheet "F20" '4.64'
Set-TextValue $newSheet "G20" '0.2227'
$newSheet.Range("H20").Value = 7

$src.Range("A2:G2").Copy($newSheet.Range("A21:G21"))
$newSheet.Range("A21").Value = 19
Set-TextValue $newSheet "B21" '001195'
$newSheet.Range("C21").Value = '工银瑞信农业产业股票'
Set-TextValue $newSheet "D21" '5.71'
Set-TextValue $newSheet "E21" '93.50'
Set-TextValue $newSheet "F21" '3.55'
Set-TextValue $newSheet "G21" '0.2027'
$newSheet.Range("H21").Value = 9

$src.Range("A2:G2").Copy($newSheet.Range("A22:G22"))
$newSheet.Range("A22").Value = 20
Set-TextValue $newSheet "B22" '900008'
$newSheet.Range("C22").Value = '中信证券稳健回报混合A'
Set-TextValue $newSheet "D22" '6.36'
Set-TextValue $newSheet "E22" '88.19'
Set-TextValue $newSheet "F22" '2.98'
Set-TextValue $newSheet "G22" '0.1895'
$newSheet.Range("H22").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A23:G23"))
$newSheet.Range("A23").Value = 21
Set-TextValue $newSheet "B23" '180020'
$newSheet.Range("C23").Value = '银华成长先锋混合'
Set-TextValue $newSheet "D23" '3.05'
Set-TextValue $newSheet "E23" '79.81'
Set-TextValue $newSheet "F23" '5.79'
Set-TextValue $newSheet "G23" '0.1766'
$newSheet.Range("H23").Value = 4

$src.Range("A2:G2").Copy($newSheet.Range("A24:G24"))
$newSheet.Range("A24").Value = 22
Set-TextValue $newSheet "B24" '671030'
$newSheet.Range("C24").Value = '西部利得事件驱动股票'
Set-TextValue $newSheet "D24" '3.55'
Set-TextValue $newSheet "E24" '94.61'
Set-TextValue $newSheet "F24" '4.97'
Set-TextValue $newSheet "G24" '0.1764'
$newSheet.Range("H24").Value = 4

$src.Range("A2:G2").Copy($newSheet.Range("A25:G25"))
$newSheet.Range("A25").Value = 23
Set-TextValue $newSheet "B25" '180018'
$newSheet.Range("C25").Value = '银华和谐主题混合'
Set-TextValue $newSheet "D25" '3.69'
Set-TextValue $newSheet "E25" '72.51'
Set-TextValue $newSheet "F25" '4.36'
Set-TextValue $newSheet "G25" '0.1609'
$newSheet.Range("H25").Value = 6

$src.Range("A2:G2").Copy($newSheet.Range("A26:G26"))
$newSheet.Range("A26").Value = 24
Set-TextValue $newSheet "B26" '014107'
$newSheet.Range("C26").Value = '博时品质生活混合A'
Set-TextValue $newSheet "D26" '4.11'
Set-TextValue $newSheet "E26" '69.12'
Set-TextValue $newSheet "F26" '3.19'
Set-TextValue $newSheet "G26" '0.1311'
$newSheet.Range("H26").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A27:G27"))
$newSheet.Range("A27").Value = 25
Set-TextValue $newSheet "B27" '008671'
$newSheet.Range("C27").Value = '银华科技创新混合'
Set-TextValue $newSheet "D27" '2.54'
Set-TextValue $newSheet "E27" '86.23'
Set-TextValue $newSheet "F27" '4.96'
Set-TextValue $newSheet "G27" '0.1260'
$newSheet.Range("H27").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A28:G28"))
$newSheet.Range("A28").Value = 26
Set-TextValue $newSheet "B28" '001277'
$newSheet.Range("C28").Value = '博时国企改革主题股票'
Set-TextValue $newSheet "D28" '2.77'
Set-TextValue $newSheet "E28" '89.22'
Set-TextValue $newSheet "F28" '3.97'
Set-TextValue $newSheet "G28" '0.1100'
$newSheet.Range("H28").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A29:G29"))
$newSheet.Range("A29").Value = 27
Set-TextValue $newSheet "B29" '009959'
$newSheet.Range("C29").Value = '长安鑫悦消费驱动混合C'
Set-TextValue $newSheet "D29" '1.62'
Set-TextValue $newSheet "E29" '86.78'
Set-TextValue $newSheet "F29" '6.15'
Set-TextValue $newSheet "G29" '0.0996'
$newSheet.Range("H29").Value = 3

$src.Range("A2:G2").Copy($newSheet.Range("A30:G30"))
$newSheet.Range("A30").Value = 28
Set-TextValue $newSheet "B30" '001940'
$newSheet.Range("C30").Value = '农银汇理现代农业加灵活配置混合'
Set-TextValue $newSheet "D30" '1.33'
Set-TextValue $newSheet "E30" '64.68'
Set-TextValue $newSheet "F30" '4.36'
Set-TextValue $newSheet "G30" '0.0580'
$newSheet.Range("H30").Value = 3

$src.Range("A2:G2").Copy($newSheet.Range("A31:G31"))
$newSheet.Range("A31").Value = 29
Set-TextValue $newSheet "B31" '900078'
$newSheet.Range("C31").Value = '中信证券稳健回报混合C'
Set-TextValue $newSheet "D31" '1.66'
Set-TextValue $newSheet "E31" '88.19'
Set-TextValue $newSheet "F31" '2.98'
Set-TextValue $newSheet "G31" '0.0495'
$newSheet.Range("H31").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A32:G32"))
$newSheet.Range("A32").Value = 30
Set-TextValue $newSheet "B32" '000826'
$newSheet.Range("C32").Value = '广发中证百度百发策略100指数A'
Set-TextValue $newSheet "D32" '4.11'
Set-TextValue $newSheet "E32" '92.42'
Set-TextValue $newSheet "F32" '1.08'
Set-TextValue $newSheet "G32" '0.0444'
$newSheet.Range("H32").Value = 6

$src.Range("A2:G2").Copy($newSheet.Range("A33:G33"))
$newSheet.Range("A33").Value = 31
Set-TextValue $newSheet "B33" '000827'
$newSheet.Range("C33").Value = '广发中证百度百发策略100指数E'
Set-TextValue $newSheet "D33" '4.11'
Set-TextValue $newSheet "E33" '92.42'
Set-TextValue $newSheet "F33" '1.08'
Set-TextValue $newSheet "G33" '0.0444'
$newSheet.Range("H33").Value = 6

$src.Range("A2:G2").Copy($newSheet.Range("A34:G34"))
$newSheet.Range("A34").Value = 32
Set-TextValue $newSheet "B34" '161913'
$newSheet.Range("C34").Value = '万家社会责任18个月定期开放混合（LOF）C'
Set-TextValue $newSheet "D34" '0.44'
Set-TextValue $newSheet "E34" '88.11'
Set-TextValue $newSheet "F34" '9.57'
Set-TextValue $newSheet "G34" '0.0421'
$newSheet.Range("H34").Value = 1

$src.Range("A2:G2").Copy($newSheet.Range("A35:G35"))
$newSheet.Range("A35").Value = 33
Set-TextValue $newSheet "B35" '900027'
$newSheet.Range("C35").Value = '中信证券信远一年持有期混合型集合资产管理计划A'
Set-TextValue $newSheet "D35" '0.71'
Set-TextValue $newSheet "E35" '75.94'
Set-TextValue $newSheet "F35" '3.49'
Set-TextValue $newSheet "G35" '0.0248'
$newSheet.Range("H35").Value = 1

$src.Range("A2:G2").Copy($newSheet.Range("A36:G36"))
$newSheet.Range("A36").Value = 34
Set-TextValue $newSheet "B36" '002319'
$newSheet.Range("C36").Value = '大成一带一路灵活配置混合'
Set-TextValue $newSheet "D36" '0.50'
Set-TextValue $newSheet "E36" '89.30'
Set-TextValue $newSheet "F36" '3.47'
Set-TextValue $newSheet "G36" '0.0174'
$newSheet.Range("H36").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A37:G37"))
$newSheet.Range("A37").Value = 35
Set-TextValue $newSheet "B37" '009855'
$newSheet.Range("C37").Value = '中加新兴成长混合A'
Set-TextValue $newSheet "D37" '0.37'
Set-TextValue $newSheet "E37" '94.73'
Set-TextValue $newSheet "F37" '2.85'
Set-TextValue $newSheet "G37" '0.0105'
$newSheet.Range("H37").Value = 10

$src.Range("A2:G2").Copy($newSheet.Range("A38:G38"))
$newSheet.Range("A38").Value = 36
Set-TextValue $newSheet "B38" '014108'
$newSheet.Range("C38").Value = '博时品质生活混合C'
Set-TextValue $newSheet "D38" '0.14'
Set-TextValue $newSheet "E38" '69.12'
Set-TextValue $newSheet "F38" '3.19'
Set-TextValue $newSheet "G38" '0.0045'
$newSheet.Range("H38").Value = 8

$src.Range("A2:G2").Copy($newSheet.Range("A39:G39"))
$newSheet.Range("A39").Value = 37
Set-TextValue $newSheet "B39" '009856'
$newSheet.Range("C39").Value = '中加新兴成长混合C'
Set-TextValue $newSheet "D39" '0.15'
Set-TextValue $newSheet "E39" '94.73'
Set-TextValue $newSheet "F39" '2.85'
Set-TextValue $newSheet "G39" '0.0043'
$newSheet.Range("H39").Value = 10

$src.Range("A2:G2").Copy($newSheet.Range("A40:G40"))
$newSheet.Range("A40").Value = 38
Set-TextValue $newSheet "B40" '900087'
$newSheet.Range("C40").Value = '中信证券信远一年持有期混合型集合资产管理计划C'
Set-TextValue $newSheet "D40" '0.02'
Set-TextValue $newSheet "E40" '75.94'
Set-TextValue $newSheet "F40" '3.49'
Set-TextValue $newSheet "G40" '0.0007'
$newSheet.Range("H40").Value = 1

$src.Range("A2:G2").Copy($newSheet.Range("A41:G41"))
$newSheet.Range("A41").Value = 39
Set-TextValue $newSheet "B41" '900077'
$newSheet.Range("C41").Value = '中信证券信远一年持有期混合型集合资产管理计划B'
Set-TextValue $newSheet "D41" '0.01'
Set-TextValue $newSheet "E41" '75.94'
Set-TextValue $newSheet "F41" '3.49'
Set-TextValue $newSheet "G41" '0.0003'
$newSheet.Range("H41").Value = 1

